# Auto-generated edit script: update cached market-price / profit values
# in the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) to match
# a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 698.34485
$ws.Range("I28").Value = 672.913
$ws.Range("J28").Value = 795.8333
$ws.Range("K28").Value = 672.913
$ws.Range("L28").Value = 795.8333
$ws.Range("M28").Value = -187.913
$ws.Range("N28").Value = -1765.8333
$ws.Range("H38").Value = 599.5833
$ws.Range("I38").Value = 102.625
$ws.Range("J38").Value = 1593.5
$ws.Range("K38").Value = 307.875
$ws.Range("L38").Value = 4780.5
$ws.Range("M38").Value = 64.125
$ws.Range("N38").Value = -5524.5
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H64").Value = 500000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 500000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 500000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -500496
$ws.Range("H67").Value = 500000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 500000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 500000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -501716
$ws.Range("H76").Value = 3268.9656
$ws.Range("I76").Value = 3268.9656
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3268.9656
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2953.9656
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3268.9656
$ws.Range("I79").Value = 3268.9656
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3268.9656
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2176.9656
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 2673.139
$ws.Range("I98").Value = 2767.4827
$ws.Range("J98").Value = 2282.2856
$ws.Range("K98").Value = 2767.4827
$ws.Range("L98").Value = 2282.2856
$ws.Range("M98").Value = -1269.4827
$ws.Range("N98").Value = -5278.2856
$ws.Range("H112").Value = 83334696
$ws.Range("J112").Value = 100001520
$ws.Range("L112").Value = 300004560
$ws.Range("N112").Value = -300006776
$ws.Range("H122").Value = 2673.139
$ws.Range("I122").Value = 2767.4827
$ws.Range("J122").Value = 2282.2856
$ws.Range("K122").Value = 8302.4481
$ws.Range("L122").Value = 6846.8568
$ws.Range("M122").Value = -5852.4481
$ws.Range("N122").Value = -11746.8568
$ws.Range("H132").Value = 349027.06
$ws.Range("I132").Value = 389037.88
$ws.Range("J132").Value = 2266.6667
$ws.Range("K132").Value = 1167113.64
$ws.Range("L132").Value = 6800.000100000001
$ws.Range("M132").Value = -1164583.64
$ws.Range("N132").Value = -11860.0001

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9159.301
$ws.Range("I32").Value = 4152.4507
$ws.Range("J32").Value = 38783.168
$ws.Range("K32").Value = 4152.4507
$ws.Range("L32").Value = 38783.168
$ws.Range("M32").Value = -3865.4507
$ws.Range("N32").Value = -39357.168
$ws.Range("H74").Value = 3941.0908
$ws.Range("I74").Value = 584.68964
$ws.Range("J74").Value = 28275
$ws.Range("K74").Value = 584.68964
$ws.Range("L74").Value = 28275
$ws.Range("M74").Value = 289.3103599999999
$ws.Range("N74").Value = -30023
$ws.Range("H77").Value = 3941.0908
$ws.Range("I77").Value = 584.68964
$ws.Range("J77").Value = 28275
$ws.Range("K77").Value = 2923.4482
$ws.Range("L77").Value = 141375
$ws.Range("M77").Value = 1444.5518
$ws.Range("N77").Value = -150111
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774
$ws.Range("H114").Value = 34598.332
$ws.Range("J114").Value = 34598.332
$ws.Range("L114").Value = 34598.332
$ws.Range("N114").Value = -43276.332
$ws.Range("H117").Value = 31700
$ws.Range("J117").Value = 31700
$ws.Range("L117").Value = 31700
$ws.Range("N117").Value = -40878
$ws.Range("H119").Value = 30500
$ws.Range("J119").Value = 30500
$ws.Range("L119").Value = 30500
$ws.Range("N119").Value = -40176
$ws.Range("H121").Value = 31000
$ws.Range("J121").Value = 31000
$ws.Range("L121").Value = 31000
$ws.Range("N121").Value = -34494
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1936.6666
$ws.Range("I105").Value = 2058.182
$ws.Range("J105").Value = 1602.5
$ws.Range("K105").Value = 2058.182
$ws.Range("L105").Value = 1602.5
$ws.Range("M105").Value = -311.1819999999998
$ws.Range("N105").Value = -5096.5

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3291
$ws.Range("I62").Value = 3082
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3082
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2458
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3291
$ws.Range("I65").Value = 3082
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 15410
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -12290
$ws.Range("N65").Value = -23740
$ws.Range("H107").Value = 1739362.6
$ws.Range("I107").Value = 2607297.5
$ws.Range("K107").Value = 2607297.5
$ws.Range("M107").Value = -2605377.5

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1219.6805
$ws.Range("I68").Value = 874.5714
$ws.Range("J68").Value = 1439.2954
$ws.Range("K68").Value = 2623.7142
$ws.Range("L68").Value = 4317.8862
$ws.Range("M68").Value = -1812.7142
$ws.Range("N68").Value = -5939.8862
$ws.Range("H71").Value = 1219.6805
$ws.Range("I71").Value = 874.5714
$ws.Range("J71").Value = 1439.2954
$ws.Range("K71").Value = 7871.1426
$ws.Range("L71").Value = 12953.6586
$ws.Range("M71").Value = -3815.1426
$ws.Range("N71").Value = -21065.6586
$ws.Range("H86").Value = 450
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 450
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356
$ws.Range("H92").Value = 828.4286
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 833.1667
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 2499.5001
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -4995.5001
$ws.Range("H107").Value = 74755.04
$ws.Range("I107").Value = 40675.28
$ws.Range("K107").Value = 122025.84
$ws.Range("M107").Value = -120105.84
$ws.Range("H131").Value = 1372377
$ws.Range("I131").Value = 7125.5557
$ws.Range("J131").Value = 1564365.5
$ws.Range("K131").Value = 21376.6671
$ws.Range("L131").Value = 4693096.5
$ws.Range("M131").Value = -16336.6671
$ws.Range("N131").Value = -4703176.5

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 21687.375
$ws.Range("J110").Value = 21687.375
$ws.Range("L110").Value = 21687.375
$ws.Range("N110").Value = -29867.375
$ws.Range("H122").Value = 7800.5625
$ws.Range("I122").Value = 12601
$ws.Range("J122").Value = 1628.5714
$ws.Range("K122").Value = 37803
$ws.Range("L122").Value = 4885.7142
$ws.Range("M122").Value = -35353
$ws.Range("N122").Value = -9785.7142

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 30640
$ws.Range("J119").Value = 30640
$ws.Range("L119").Value = 30640
$ws.Range("N119").Value = -40316
$ws.Range("H132").Value = 4894.4346
$ws.Range("I132").Value = 5106.6943
$ws.Range("J132").Value = 4130.3
$ws.Range("K132").Value = 15320.0829
$ws.Range("L132").Value = 12390.9
$ws.Range("M132").Value = -12790.0829
$ws.Range("N132").Value = -17450.9
$ws.Range("H136").Value = 982.8125
$ws.Range("I136").Value = 1025.2858
$ws.Range("J136").Value = 901.7273
$ws.Range("K136").Value = 3075.8574
$ws.Range("L136").Value = 2705.1819
$ws.Range("M136").Value = -525.8574000000003
$ws.Range("N136").Value = -7805.1819

